$d = $word.ActiveDocument

$values = @(
    "12+8=", "5+85=", "61-3=", "40+32=", "79-38=",
    "34-0=", "22+65=", "61-18=", "46-17=", "30+50=",
    "9+83=", "0+95=", "23+55=", "76-62=", "90-19=",
    "97+2=", "41+0=", "17-0=", "74+18=", "4+70=",
    "90-54=", "90-9=", "96-32=", "20+49=", "45+22=",
    "48+21=", "99-64=", "19+69=", "31+68=", "40+23=",
    "9+6=", "93-40=", "92+7=", "98-81=", "25+52=",
    "23+50=", "71+24=", "17+3=", "58-32=", "94-62=",
    "72+10=", "29+26=", "30+1=", "87-30=", "64-20=",
    "14-4=", "61-36=", "44+25=", "50+26=", "43+20=",
    "51-48=", "27+40=", "60-28=", "36+60=", "97-14=",
    "74+0=", "3+45=", "22+13=", "74-56=", "36-0=",
    "67-26=", "97-4=", "35+32=", "8+31=", "26+13=",
    "81-3=", "4+2=", "47-26=", "78-20=", "69-28=",
    "69-4=", "73-60=", "19-18=", "71-5=", "98-37=",
    "24+12=", "34-21=", "68+14=", "4+74=", "57+30=",
    "36-18=", "23+70=", "40+22=", "17-6=", "54+22=",
    "33-6=", "76+19=", "24-21=", "52-39=", "44+51=",
    "1+94=", "71-35=", "42-41=", "18-14=", "38+0=",
    "23+25=", "36-23=", "91-74=", "99-11=", "53-53="
)

$t = $d.Tables.Item(1)
$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $cell = $row.Cells.Item($c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
